$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "bdvb(12)"
$ws.Range("A3").Value = "gggggg(gg)"
$ws.Range("A4").Value = "llll(ll)"
$ws.Range("A5").Value = "Xiao(xiao)"
